# Super Gnomatic Rock Cycle - add the "Extrusive" igneous rocks (pumice,
# scoria, rhyolite, basalt) as new Key/Value/VoiceDuration rows, matching
# the existing sedimentary / other igneous rock-entry layout, then move
# the on-screen selection/goal marker down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- pumice -----------------------------------------------------------
$ws.Range("A20").Value = "pumice"
$ws.Range("B20").Value = "Pumice"
$ws.Range("C20").Value = 0.6

$ws.Range("A21").Value = "pumiceDesc"
$ws.Range("B21").Value = "Igneous rock. Extrusive"
$ws.Range("C21").Value = 5

# --- scoria -------------------------------------------------------------
$ws.Range("A22").Value = "scoria"
$ws.Range("B22").Value = "Scoria"
$ws.Range("C22").Value = 0.6

$ws.Range("A23").Value = "scoriaDesc"
$ws.Range("B23").Value = "Igneous rock. Extrusive"
$ws.Range("C23").Value = 5

# --- rhyolite (value cell typed before the key cell, same as source) ---
$ws.Range("B24").Value = "Rhyolite"
$ws.Range("A24").Value = "rhyolite"
$ws.Range("C24").Value = 0.6

$ws.Range("A25").Value = "rhyoliteDesc"
$ws.Range("B25").Value = "Igneous rock. Extrusive"
$ws.Range("C25").Value = 5

# --- basalt ---------------------------------------------------------
$ws.Range("A26").Value = "basalt"
$ws.Range("B26").Value = "Basalt"
$ws.Range("C26").Value = 0.6

$ws.Range("A27").Value = "basaltDesc"
$ws.Range("B27").Value = "Igneous rock. Extrusive"
$ws.Range("C27").Value = 5

# --- offscreen goal display: move selection to the new bottom row ------
$ws.Activate()
$ws.Range("A26").Select()
